# Customer data correction: the customer previously recorded in row 3
# ("Alok Patel") is actually "Anjana Singh". Update the Customer Name
# cell for that row accordingly (mirrors the data fix captured in the
# commit's refactored test-data generation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "Anjana Singh"
